$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "both" actionType option was dropped; rows that used it now use "contact" instead.
$ws.Range("H2").Value = "contact"
$ws.Range("H3").Value = "contact"

# Reflect the new selection left behind in the sheet after the edit.
[void]$ws.Range("J14").Select()
